$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 118338939.2758621
$ws.Range("D2").Value = 2.487842556667488
$ws.Range("E2").Value = 1.260517790344322
$ws.Range("F2").Value = 1.260517790344322
$ws.Range("C3").Value = 117723824.8666667
$ws.Range("D3").Value = 2.49074680761321
$ws.Range("E3").Value = 1.265626611507313
$ws.Range("F3").Value = 2.526144401851635
$ws.Range("C4").Value = 118075866.27
$ws.Range("D4").Value = 2.387410944463444
$ws.Range("E4").Value = 32.28313493547573
$ws.Range("F4").Value = 34.80927933732737
$ws.Range("C5").Value = 118091567.1275
$ws.Range("D5").Value = 2.381252895741406
$ws.Range("E5").Value = 32.36231756547575
$ws.Range("F5").Value = 67.17159690280312
$ws.Range("C6").Value = 118166435.175
$ws.Range("D6").Value = 2.365426027882202
$ws.Range("E6").Value = 13.74173268829829
$ws.Range("F6").Value = 80.9133295911014
$ws.Range("C7").Value = 118066228.4025
$ws.Range("D7").Value = 2.294032299898257
$ws.Range("E7").Value = 14.48063077259742
$ws.Range("F7").Value = 95.39396036369882
$ws.Range("C8").Value = 118033938.22
$ws.Range("D8").Value = 2.300343855353062
$ws.Range("E8").Value = 14.37089871702059
$ws.Range("F8").Value = 109.7648590807194
$ws.Range("C9").Value = 118070018.7425
$ws.Range("D9").Value = 2.078483567790478
$ws.Range("E9").Value = 83.53483776515513
$ws.Range("F9").Value = 193.2996968458745
$ws.Range("C10").Value = 118075682.8225
$ws.Range("D10").Value = 2.1169956413529
$ws.Range("E10").Value = 82.41130489106892
$ws.Range("F10").Value = 275.7110017369434
$ws.Range("C11").Value = 117977392.25
$ws.Range("D11").Value = 2.068225692092291
$ws.Range("E11").Value = 84.42488918090511
$ws.Range("F11").Value = 360.1358909178485
$ws.Range("C12").Value = 118148119.3325
$ws.Range("D12").Value = 2.066854092554542
$ws.Range("E12").Value = 22.67535315631811
$ws.Range("F12").Value = 382.8112440741667
$ws.Range("C13").Value = 118106334.5425
$ws.Range("D13").Value = 2.09851492523725
$ws.Range("E13").Value = 22.66812028049555
$ws.Range("F13").Value = 405.4793643546622
$ws.Range("C14").Value = 118076243.5325
$ws.Range("D14").Value = 2.105226826843285
$ws.Range("E14").Value = 22.5208315152537
$ws.Range("F14").Value = 428.0001958699159
$ws.Range("C15").Value = 118072455.6907731
$ws.Range("D15").Value = 2.279505060492302
$ws.Range("E15").Value = 50.53005277593896
$ws.Range("F15").Value = 478.5302486458548
$ws.Range("C16").Value = 118033307.3875
$ws.Range("D16").Value = 2.295387219711106
$ws.Range("E16").Value = 50.56616661698777
$ws.Range("F16").Value = 529.0964152628426
$ws.Range("C17").Value = 118125114.985
$ws.Range("D17").Value = 2.276925417356452
$ws.Range("E17").Value = 49.44935048586963
$ws.Range("F17").Value = 578.5457657487123
$ws.Range("C18").Value = 118090900.7775
$ws.Range("D18").Value = 2.227327930608984
$ws.Range("E18").Value = 14.59319302190979
$ws.Range("F18").Value = 593.138958770622
$ws.Range("C19").Value = 118096560.6225
$ws.Range("D19").Value = 2.082177426834825
$ws.Range("E19").Value = 15.93931734344292
$ws.Range("F19").Value = 609.0782761140649
$ws.Range("C20").Value = 118096630.625
$ws.Range("D20").Value = 2.073989146208124
$ws.Range("E20").Value = 15.92025863796721
$ws.Range("F20").Value = 624.9985347520321
